$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.930923581123352
$ws.Range("B1").Value = 3.262081623077393
$ws.Range("C1").Value = 2.823889255523682
$ws.Range("D1").Value = 0.9937343001365662
$ws.Range("E1").Value = 0.6492149829864502
